$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the new "Save" header value in H1
$ws.Range("H1").Value = "Save"

# Copy formatting from the neighboring header cell (G1) onto H1
# so the new column matches the existing header style (bold, centered, bordered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the numeric value for the new "Save" column in row 2
$ws.Range("H2").Value = 0
